$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player names for rows 23-41 (score column stays 200 except row31->1200, row38 reverts to 200)
$ws.Range("A23").Value = "NIKOLA"
$ws.Range("A24").Value = "RIMA"
$ws.Range("A25").Value = "OMARAS"
$ws.Range("A26").Value = "PATRICIJA"
$ws.Range("A27").Value = "GERALDAS"
$ws.Range("A28").Value = "MARINA"
$ws.Range("A29").Value = "MARIJA"
$ws.Range("A30").Value = "DOMAS"
$ws.Range("A31").Value = "DENISAS"
$ws.Range("B31").Value = 1200
$ws.Range("A32").Value = "BOBAS"
$ws.Range("A33").Value = "UGNE"
$ws.Range("A34").Value = "VITALIJUS"
$ws.Range("A35").Value = "LAJA"
$ws.Range("A36").Value = "ROMAS"
$ws.Range("A37").Value = "RIMA"
$ws.Range("A38").Value = "ERNESTAS"
$ws.Range("B38").Value = 200
$ws.Range("A39").Value = "ONA"
$ws.Range("A40").Value = "DOMANTAS"
$ws.Range("B40").Value = 200
$ws.Range("A41").Value = "ULA"
$ws.Range("B41").Value = 200

# Remove old trailing rows 42-44 (table now ends at row 41)
$ws.Range("A42:B44").ClearContents()

# Autofit column A width
$ws.Columns.Item(1).AutoFit()

# Update selection to match new active cell
$ws.Range("C8").Select()
